# Insert one new data row at row 205 (pushes the existing rows 205..289
# down to 206..290, growing the used range from A1:R289 to A1:R290), then
# populate the new row with a new weekly Espinaca price observation for
# "Terminal La Palmera de La Serena".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(205).Insert()

$ws.Range("A205").Value = 8
$ws.Range("B205").Value = "Terminal La Palmera de La Serena"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = [DateTime]::FromOADate(44755)
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = 100112012
$ws.Range("G205").Value = "Espinaca"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 2800
$ws.Range("K205").Value = 500
$ws.Range("L205").Value = 600
$ws.Range("M205").Value = 550
$ws.Range("N205").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O205").Value = "Provincia del Elquí"
$ws.Range("P205").Value = 1100
$ws.Range("Q205").Value = 0.5
$ws.Range("R205").Value = "Hortaliza"
